$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.402.02"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.51"
$ws.Range("E3").Value = "  -1.54%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.09"
$ws.Range("E5").Value = "  -2.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.532"
$ws.Range("E6").Value = "  +3.50%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.14"
$ws.Range("E8").Value = "  -1.66%  "

$ws.Range("E9").Value = "  -3.27%  "

$ws.Range("E10").Value = "  -2.22%  "

$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.93"
$ws.Range("E12").Value = "  -1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.637.66"
$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  -2.73%  "

$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.21"
$ws.Range("E16").Value = "  -2.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.391.49"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.32"
$ws.Range("E18").Value = "  -5.35%  "

$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("E22").Value = "  -3.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("E24").Value = "  -1.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.69"
$ws.Range("E25").Value = "  +1.11%  "

$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.94"
$ws.Range("E27").Value = "  -3.46%  "

$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.51"
$ws.Range("E29").Value = "  -5.12%  "

$ws.Range("E30").Value = "  -4.96%  "

$ws.Range("E31").Value = "  -3.89%  "

$ws.Range("E32").Value = "  -1.96%  "

$ws.Range("E33").Value = "  -0.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.409.93"
$ws.Range("E34").Value = "  -4.60%  "

$ws.Range("E35").Value = "  -0.20%  "

$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.561"
$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("E38").Value = "  -5.89%  "

$ws.Range("E39").Value = "  -2.98%  "

$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.47"
$ws.Range("E42").Value = "  -1.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("E44").Value = "  +0.96%  "

$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.48"
$ws.Range("E46").Value = "  -7.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.782.74"
$ws.Range("E47").Value = "  -1.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.63"
$ws.Range("E48").Value = "  -4.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.18"
$ws.Range("E49").Value = "  -2.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0104"
$ws.Range("E50").Value = "  -3.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0987"
$ws.Range("E51").Value = "  -3.71%  "
